$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 7, car formerly "obs..." id 19) entirely.
$ws.Rows(7).Delete()

# Clear cells whose columns are no longer populated for each remaining row.
$ws.Range("M2").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("K6").ClearContents()

# Row 2 (car_01)
$ws.Range("A2").Value = "car_01"
$ws.Range("B2").Value = 11.00746188820342
$ws.Range("C2").Value = 45.43949340167618
$ws.Range("D2").Value = 308
$ws.Range("E2").Value = "(car_02, 110.05), (car_03, 122.93), (car_04, 67.75), (car_05, 75.01), (car_02, 110.05), (car_03, 122.93), (car_04, 67.75), (car_05, 75.01)"
$ws.Range("F2").Value = "(obs01, 38.8), (obs02, 5.18), (obs03, 5.71), (obs04, 4.83), (obs05, 24.92), (obs06, 45.9), (obs08, 34.96), (obs09, 2.16), (obs10, 5.02), (obs01, 38.8), (obs02, 5.18), (obs03, 5.71), (obs04, 4.83), (obs05, 24.92), (obs06, 45.9), (obs08, 34.96), (obs09, 2.16), (obs10, 5.02)"
$ws.Range("G2").Value = "(obs01, 38.8), (obs01, 38.8)"
$ws.Range("H2").Value = "(obs05, 24.92), (obs08, 34.96), (obs05, 24.92), (obs08, 34.96)"
$ws.Range("I2").Value = "(obs03, 5.71), (obs04, 4.83), (obs10, 5.02), (obs03, 5.71), (obs04, 4.83), (obs10, 5.02)"
$ws.Range("J2").Value = "(obs02, 5.18), (obs06, 45.9), (obs09, 2.16), (obs02, 5.18), (obs06, 45.9), (obs09, 2.16)"
$ws.Range("N2").Value = "(car_02, 110.05), (car_03, 122.93), (car_04, 67.75), (car_05, 75.01), (car_02, 110.05), (car_03, 122.93), (car_04, 67.75), (car_05, 75.01)"

# Row 3 (car_02)
$ws.Range("A3").Value = "car_02"
$ws.Range("B3").Value = 11.00844905524499
$ws.Range("C3").Value = 45.43961881483433
$ws.Range("D3").Value = 101
$ws.Range("E3").Value = "(car_01, 110.05), (car_03, 13.09), (car_04, 42.49), (car_05, 35.09), (car_01, 110.05), (car_03, 13.09), (car_04, 42.49), (car_05, 35.09)"
$ws.Range("L3").Value = "(car_03, 13.09), (car_03, 13.09)"
$ws.Range("N3").Value = "(car_01, 110.05), (car_04, 42.49), (car_05, 35.09), (car_01, 110.05), (car_04, 42.49), (car_05, 35.09)"

# Row 4 (car_03)
$ws.Range("A4").Value = "car_03"
$ws.Range("B4").Value = 11.00856714629746
$ws.Range("C4").Value = 45.43961088955921
$ws.Range("D4").Value = 134
$ws.Range("E4").Value = "(car_01, 122.93), (car_02, 13.09), (car_04, 55.51), (car_05, 47.92), (car_01, 122.93), (car_02, 13.09), (car_04, 55.51), (car_05, 47.92)"
$ws.Range("K4").Value = "(car_02, 13.09), (car_02, 13.09)"
$ws.Range("N4").Value = "(car_01, 122.93), (car_04, 55.51), (car_05, 47.92), (car_01, 122.93), (car_04, 55.51), (car_05, 47.92)"

# Row 5 (car_04)
$ws.Range("A5").Value = "car_04"
$ws.Range("B5").Value = 11.00806546558897
$ws.Range("C5").Value = 45.43959899859149
$ws.Range("D5").Value = 307
$ws.Range("E5").Value = "(car_01, 67.75), (car_02, 42.49), (car_03, 55.51), (car_05, 8.67), (car_01, 67.75), (car_02, 42.49), (car_03, 55.51), (car_05, 8.67)"
$ws.Range("F5").Value = "(obs01, 41.95), (obs06, 22.19), (obs07, 14.0), (obs01, 41.95), (obs06, 22.19), (obs07, 14.0)"
$ws.Range("H5").Value = "(obs01, 41.95), (obs06, 22.19), (obs07, 14.0), (obs01, 41.95), (obs06, 22.19), (obs07, 14.0)"
$ws.Range("K5").Value = "(car_05, 8.67), (car_05, 8.67)"
$ws.Range("L5").Value = "(car_01, 67.75), (car_01, 67.75)"
$ws.Range("N5").Value = "(car_02, 42.49), (car_03, 55.51), (car_02, 42.49), (car_03, 55.51)"

# Row 6 (car_05)
$ws.Range("A6").Value = "car_05"
$ws.Range("B6").Value = 11.00813626526003
$ws.Range("C6").Value = 45.43956504282627
$ws.Range("D6").Value = 278
$ws.Range("E6").Value = "(car_01, 75.01), (car_02, 35.09), (car_03, 47.92), (car_04, 8.67), (car_01, 75.01), (car_02, 35.09), (car_03, 47.92), (car_04, 8.67)"
$ws.Range("F6").Value = "(obs01, 46.07), (obs06, 30.16), (obs07, 21.65), (obs01, 46.07), (obs06, 30.16), (obs07, 21.65)"
$ws.Range("H6").Value = "(obs01, 46.07), (obs06, 30.16), (obs07, 21.65), (obs01, 46.07), (obs06, 30.16), (obs07, 21.65)"
$ws.Range("L6").Value = "(car_01, 75.01), (car_04, 8.67), (car_01, 75.01), (car_04, 8.67)"
$ws.Range("N6").Value = "(car_02, 35.09), (car_03, 47.92), (car_02, 35.09), (car_03, 47.92)"
